# Generate Report for Handback
# This mirrors a localization-status report refresh: the status text moves
# from "Ready for handoff" to "Handed back: in sync with en-US", and the two
# rows on each language sheet get their handback file + handback datetime
# filled in (with a hyperlink on the "Latest Target File" column).

$wb = $excel.ActiveWorkbook

$ws_overview = $wb.Worksheets.Item("Overview")
$ws_zhcn     = $wb.Worksheets.Item("zh-cn")
$ws_dede     = $wb.Worksheets.Item("de-de")

$newStatus = "Handed back: in sync with en-US"

# --- Overview sheet: zh-cn / de-de status columns (E, F) for both rows ---
$ws_overview.Range("E2").Value = $newStatus
$ws_overview.Range("F2").Value = $newStatus
$ws_overview.Range("E3").Value = $newStatus
$ws_overview.Range("F3").Value = $newStatus

$aUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/66114aee9abbaade7769563da60cbea8375cbbcc/e2e/a.md"

# --- zh-cn sheet ---
$ws_zhcn.Range("C2").Value = $newStatus
$ws_zhcn.Range("C3").Value = $newStatus

$ws_zhcn.Range("I2").Value = "a.md"
$ws_zhcn.Range("J2").Value = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf"
$ws_zhcn.Range("K2").Value = "2016-08-24 20:38:19"

$ws_zhcn.Range("I3").Value = "a.md"
$ws_zhcn.Range("J3").Value = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf"
$ws_zhcn.Range("K3").Value = "2016-08-24 20:38:19"

$ws_zhcn.Hyperlinks.Add($ws_zhcn.Range("I2"), $aUrl, [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, "a.md")
$ws_zhcn.Hyperlinks.Add($ws_zhcn.Range("I3"), $aUrl, [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, "a.md")

# --- de-de sheet ---
$ws_dede.Range("C2").Value = $newStatus
$ws_dede.Range("C3").Value = $newStatus

$ws_dede.Range("I2").Value = "a.md"
$ws_dede.Range("J2").Value = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf"
$ws_dede.Range("K2").Value = "2016-08-24 20:38:27"

$ws_dede.Range("I3").Value = "a.md"
$ws_dede.Range("J3").Value = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf"
$ws_dede.Range("K3").Value = "2016-08-24 20:38:27"

$ws_dede.Hyperlinks.Add($ws_dede.Range("I2"), $aUrl, [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, "a.md")
$ws_dede.Hyperlinks.Add($ws_dede.Range("I3"), $aUrl, [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, "a.md")

# --- Column widths: widen "Status"-like columns to fit the longer text, and
# widen the "Latest Handback File" column to fit the new xlf filename ---
$ws_overview.Columns.Item(5).ColumnWidth = 29.9777047293527
$ws_overview.Columns.Item(6).ColumnWidth = 29.9777047293527

$ws_zhcn.Columns.Item(3).ColumnWidth = 29.9777047293527
$ws_zhcn.Columns.Item(10).ColumnWidth = 40

$ws_dede.Columns.Item(3).ColumnWidth = 29.9777047293527
$ws_dede.Columns.Item(10).ColumnWidth = 40
